$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "Insurance"
$ws.Range("C4").Value = "Telecommunications Services"
$ws.Range("C8").Value = "Food & Tobacco"
$ws.Range("C20").Value = "Telecommunications Services"
$ws.Range("C23").Value = "Multiline Utilities"
$ws.Range("C26").Value = "Health Care Equipment & Supplies"
$ws.Range("C37").Value = "Retailers"
$ws.Range("C42").Value = "Electronic Equipment & Parts"
$ws.Range("C44").Value = "Pharmaceuticals & Biotechnology"
$ws.Range("C45").Value = "Homebuilding & Construction Supplies"
$ws.Range("C48").Value = "Machinery, Tools, Heavy Vehicles, Trains & Ships"
$ws.Range("C59").Value = "Financial Services"
$ws.Range("C60").Value = "Food & Drug Retailing"
$ws.Range("C61").Value = "Aerospace & Defence"
$ws.Range("C62").Value = "Containers & Packaging"
$ws.Range("C63").Value = "Multiline Utilities"
$ws.Range("C79").Value = "Food & Drug Retailing"
$ws.Range("C81").Value = "Collective Investments"
$ws.Range("C83").Value = "Multiline Utilities"
$ws.Range("C87").Value = "Health Care Equipment & Supplies"
$ws.Range("C90").Value = "Electrical Utilities & Independent Power Producers"
$ws.Range("C94").Value = "Food & Drug Retailing"
$ws.Range("C96").Value = "Multiline Utilities"
$ws.Range("C97").Value = "Real Estate Investment Trusts"
